# Apply weekly update of Fruta/hortaliza prices: rows 2-13 (except row 10)
# get their D,L,M,N,O,P,Q,R,S,T values reshuffled according to the source
# repository update (each row picks up the data previously held by another
# row in the same block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move, keyed by original
# row number, so that we can safely overwrite cells in place without losing
# data we still need to read later.
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$srcRows = 2..13
$snapshot = @{}
foreach ($r in $srcRows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: target row -> row whose snapshot values should be written there.
$mapping = @{
    2  = 6
    3  = 9
    4  = 13
    5  = 11
    6  = 4
    7  = 5
    8  = 2
    9  = 3
    10 = 10
    11 = 12
    12 = 7
    13 = 8
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $rowData = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $rowData[$c]
    }
}
